$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (price) cells as Text so numeric-looking strings
# (e.g. "318.80") are not auto-converted to real numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.677.98"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.474.55"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "318.80"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("D6").Value = "92.89"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.517"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("D10").Value = "33.13"
$ws.Range("E10").Value = "  +2.46%  "
$ws.Range("D11").Value = "0.0855"
$ws.Range("E11").Value = "  +8.31%  "
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").Value = "2.855.75"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "6.90"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "15.79"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").Value = "2.470.60"
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").Value = "0.788"
$ws.Range("E17").Value = "  +2.83%  "
$ws.Range("D18").Value = "41.638.39"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D20").Value = "0.0₃0951"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").Value = "11.30"
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("D23").Value = "239.54"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "24.79"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("D29").Value = "9.84"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").Value = "36.09"
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("D31").Value = "158.82"
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "2.59"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0767"
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").Value = "17.35"
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("E37").Value = "  +4.90%  "
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("E42").Value = "  +5.23%  "
$ws.Range("D43").Value = "1.996.24"
$ws.Range("E43").Value = "  +2.20%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0285"
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "18.96"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("D46").Value = "2.98"
$ws.Range("E46").Value = "  +2.18%  "
$ws.Range("E47").Value = "  +3.08%  "
$ws.Range("D48").Value = "2.712.50"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").Value = "97.30"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("E50").Value = "  +3.42%  "
$ws.Range("D51").Value = "67.12"
$ws.Range("E51").Value = "  +0.49%  "

# Restore default styling on column D so no stray number-format style
# is left attached to the cells.
$ws.Range("D2:D51").Style = "Normal"
